{"js": "// The document contains one table of two-digit x two-digit multiplication\n// problems (\"NN\u00d7NN=RRRR\"). This edit swaps the problem/answer text of every\n// populated cell (in document/reading order) for a freshly generated one,\n// while leaving every empty \"spacer\" row and all cell/paragraph/run\n// formatting untouched.\nconst newTexts = [\n  \"43\u00d785=3655\", \"21\u00d794=1974\", \"56\u00d741=2296\", \"34\u00d758=1972\", \"37\u00d745=1665\",\n  \"34\u00d753=1802\", \"32\u00d783=2656\", \"28\u00d793=2604\", \"63\u00d785=5355\", \"17\u00d715=255\",\n  \"73\u00d774=5402\", \"59\u00d769=4071\", \"81\u00d793=7533\", \"37\u00d711=407\",  \"55\u00d738=2090\",\n  \"90\u00d742=3780\", \"37\u00d744=1628\", \"86\u00d785=7310\", \"50\u00d775=3750\", \"62\u00d783=5146\",\n  \"25\u00d749=1225\", \"83\u00d715=1245\", \"14\u00d727=378\",  \"50\u00d773=3650\", \"99\u00d769=6831\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  return \"no tables found\";\n}\nconst table = tables.items[0];\n\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items/value\");\n}\nawait context.sync();\n\n// Walk the table in document order and replace the text of every\n// non-empty cell with the corresponding new value, leaving blank rows\n// (and any cell formatting) exactly as they were.\nlet idx = 0;\nfor (const row of rows) {\n  for (const cell of row.cells.items) {\n    const current = (cell.value || \"\").trim();\n    if (current.length > 0 && idx < newTexts.length) {\n      const range = cell.body.getRange(\"Whole\");\n      range.insertText(newTexts[idx], \"Replace\");\n      idx++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains one table of two-digit x two-digit multiplication\n# problems (\"NN\u00d7NN=RRRR\"). This edit swaps the problem/answer text of every\n# populated cell (in row-major / reading order) for a freshly generated\n# one, while leaving every empty \"spacer\" row and all cell/paragraph/run\n# formatting untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newTexts = @(\n  \"43\u00d785=3655\", \"21\u00d794=1974\", \"56\u00d741=2296\", \"34\u00d758=1972\", \"37\u00d745=1665\",\n  \"34\u00d753=1802\", \"32\u00d783=2656\", \"28\u00d793=2604\", \"63\u00d785=5355\", \"17\u00d715=255\",\n  \"73\u00d774=5402\", \"59\u00d769=4071\", \"81\u00d793=7533\", \"37\u00d711=407\",  \"55\u00d738=2090\",\n  \"90\u00d742=3780\", \"37\u00d744=1628\", \"86\u00d785=7310\", \"50\u00d775=3750\", \"62\u00d783=5146\",\n  \"25\u00d749=1225\", \"83\u00d715=1245\", \"14\u00d727=378\",  \"50\u00d773=3650\", \"99\u00d769=6831\"\n)\n\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    # Cell.Range.Text carries a trailing cell-mark (CR + BEL); strip it\n    # before deciding whether the cell actually holds text.\n    $clean = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($clean.Length -gt 0 -and $idx -lt $newTexts.Length) {\n      $cell.Range.Text = $newTexts[$idx]\n      $idx++\n    }\n  }\n}\n\nWrite-Output (\"replaced=\" + $idx)\n"}
